$p = $ppt.ActivePresentation

# --- Slide 2: fix "unemployeed" typo and merge the split runs into one run ---
$s2 = $p.Slides.Item(2)
$shape2 = $s2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange
$tr2.Delete()
$tr2.Text = "With many being unemployed due to Covid-19, hunger among US children is growing. Food insecurity has also grown from students not receiving a meal at schools are closed or operating virtually due to Covid-19.In April 41% of mothers with children under the age of 13 reported recent food insecurity, which is the highest level since 2001. In addition, the non-profit organization, Feeding America, projects 18 million children could be food insecure in 2020. While many school districts are operating in a grab-n-go fashion, many students are unable to make it to these sites. For example, many school districts operate summer feeding programs. However, in 2018 only 14.1% of kids who received a free or reduced-price meals during the school year got a meal over the summer. Barriers for students getting grab-n-go meals include transportation barriers, families not knowing meals are available, and sites not being open when families are able to pick up meals.`nThis project will look at where grab-n-go sites should be placed within the Lewisville ISD school district. Selecting the right grab-n-go sites is crucial in helping students overcome the barriers preventing them from getting the food they need. "

# --- Slide 3: reposition/resize the picture ---
# Target EMU: off x=2771775 y=2461463, ext cx=6429376 cy=1860415 (1 pt = 12700 EMU)
$s3 = $p.Slides.Item(3)
$pic = $s3.Shapes.Item(3)
$pic.Left = 218.25
$pic.Top = 193.8159866519685
$pic.Width = 506.2500787401575
$pic.Height = 146.48937227874018

# --- Slide 8: fix "way be"/"determing" typos and merge the split runs into one run ---
$s8 = $p.Slides.Item(8)
$shape8 = $s8.Shapes.Item(2)
$tr8 = $shape8.TextFrame.TextRange
$tr8.Delete()
$tr8.Text = "When selecting sites, stakeholders must look at how many sites the school is looking to operate and would need to look at how accessible each site is to its students. For example, parents may be working and students may have to walk to the sites. The number of sites selected will play a big role in determining the sites. If few sites are selected, it would be advantageous to select sites that are geographically separated compared to just selecting the sites with the highest need. Also, it is important to evaluate other resources that are available outside of the school district's feeding program to prevent two feeding programs from targeting the same students."
